$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 250, shifting existing rows 250:316 down to 251:317.
$ws.Rows.Item(250).Insert()

# Fill the newly inserted row 250 with the new record.
$ws.Cells.Item(250, 1).Value = 10
$ws.Cells.Item(250, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(250, 3).Value = "La Araucanía"
$ws.Cells.Item(250, 4).Value = 44736
$ws.Cells.Item(250, 5).Value = 9
$ws.Cells.Item(250, 6).Value = 100112017
$ws.Cells.Item(250, 7).Value = "Apio"
$ws.Cells.Item(250, 8).Value = "Americana (o)"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 130
$ws.Cells.Item(250, 11).Value = 8000
$ws.Cells.Item(250, 12).Value = 9000
$ws.Cells.Item(250, 13).Value = 8385
$ws.Cells.Item(250, 14).Value = "$/docena de matas"
$ws.Cells.Item(250, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(250, 16).Value = 1398
$ws.Cells.Item(250, 17).Value = 6
$ws.Cells.Item(250, 18).Value = "Hortaliza"
